$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Keyboard" (sheet2): add a new "SPACE" keyword row (8), copying the
# format from the row above it (7).
# ---------------------------------------------------------------------------
$wsKeyboard = $wb.Worksheets.Item("Keyboard")

$wsKeyboard.Range("A7:C7").Copy()
$wsKeyboard.Range("A8:C8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsKeyboard.Range("A8").Value = "SPACE"
$wsKeyboard.Range("B8").Value = "SPACE"
$wsKeyboard.Range("C8").Value = "descr."

# ---------------------------------------------------------------------------
# Sheet "Credentials" (sheet1): add two new credential rows (34 & 35) for the
# new Product Factory Neo4j database account, re-using the formatting of the
# two rows directly above (32 & 33) the same way a user would by typing into
# the row right after an existing, already-styled block. Tags (column A) are
# filled in first for both rows, then the values (column B) -- matching how
# the sheet's shared strings ended up ordered.
# ---------------------------------------------------------------------------
$wsCred = $wb.Worksheets.Item("Credentials")

$wsCred.Range("A32:C33").Copy()
$wsCred.Range("A34:C35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsCred.Range("A34").Value = "PRODUCTFACTORYDATABASEUSERNEW"
$wsCred.Range("A35").Value = "PRODUCTFACTORYDATABASEPASSWORDNEW"

$wsCred.Range("B34").Value = "Neo4j"
$wsCred.Range("B35").Value = "3M#$,qns5uw*W#jr"

$wsCred.Range("C34").Value = "descr."
$wsCred.Range("C35").Value = "descr."

# ---------------------------------------------------------------------------
# View-state changes: move the active tab/selection around, matching the
# recorded session. Order matters -- the last worksheet activated ends up
# "tabSelected" and drives the workbook's active tab.
# ---------------------------------------------------------------------------
$wsLinks = $wb.Worksheets.Item("Links")
$wsValiddata = $wb.Worksheets.Item("Validdata")

$wsKeyboard.Activate()
$excel.ActiveWindow.Zoom = 131
$wsKeyboard.Range("C7:C8").Select()

$wsLinks.Activate()
$wsLinks.Range("B20").Select()

$wsCred.Activate()
$wsCred.Range("C18").Select()
